$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.657.92'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +6.49%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.635.07'
$ws.Range('D3').ClearFormats()
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '513.92'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.78'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.00%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.613'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.995'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.676.06'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +10.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.29'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +11.62%  '
$ws.Range('E11').Value = '  +5.84%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.37%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.100.33'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +9.29%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '60.821.97'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.57%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.88'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +5.52%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000140'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.48%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.667.92'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +9.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.80'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '351.04'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +8.21%  '
$ws.Range('E21').Value = '  +5.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.19'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.28%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.76%  '
$ws.Range('E25').Value = '  +3.65%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.771.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +9.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.166'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +5.16%  '
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0871'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +11.54%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.55'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.60%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.62'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '157.55'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.11%  '
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.74'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +8.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.04'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.38%  '
$ws.Range('E37').Value = '  +6.42%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.52'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.60%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.871'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '310.94'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +17.61%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.78'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +7.67%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.834'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +29.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '35.56'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.25%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.647'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +9.35%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0578'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +8.89%  '
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '20.08'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +15.32%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.04'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.23%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.991'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  +4.38%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.039.62'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +9.50%  '
